$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 19).Value = 1.31  # S2: 1.3 -> 1.31
$ws.Cells.Item(4, 8).Value = 3.35  # H4: 3.4 -> 3.35
$ws.Cells.Item(4, 9).Value = 5.8  # I4: 5.7 -> 5.8
$ws.Cells.Item(4, 17).Value = 2.72  # Q4: 1.01 -> 2.72
$ws.Cells.Item(4, 18).Value = 1.13  # R4: 1.12 -> 1.13
$ws.Cells.Item(4, 19).Value = 2.72  # S4: 1.01 -> 2.72
$ws.Cells.Item(5, 6).Value = 1.34  # F5: 1.33 -> 1.34
$ws.Cells.Item(5, 7).Value = 1.59  # G5: 1.61 -> 1.59
$ws.Cells.Item(5, 8).Value = 1.04  # H5: 8.6 -> 1.04
$ws.Cells.Item(5, 10).Value = 3.7  # J5: 3.6 -> 3.7
$ws.Cells.Item(5, 11).Value = 980  # K5: 7.6 -> 980
$ws.Cells.Item(5, 12).Value = 1.4  # L5: 1.01 -> 1.4
$ws.Cells.Item(5, 13).Value = 1.08  # M5: 1.01 -> 1.08
$ws.Cells.Item(5, 14).Value = 1.01  # N5: 1.51 -> 1.01
$ws.Cells.Item(5, 15).Value = 1.31  # O5: 1.01 -> 1.31
$ws.Cells.Item(5, 16).Value = 1.52  # P5: 1.51 -> 1.52
$ws.Cells.Item(5, 17).Value = 2.08  # Q5: 2.04 -> 2.08
$ws.Cells.Item(5, 18).Value = 1.16  # R5: 1.13 -> 1.16
$ws.Cells.Item(5, 19).Value = 3.45  # S5: 2.04 -> 3.45
$ws.Cells.Item(5, 20).Value = 2.56  # T5: 1.01 -> 2.56
$ws.Cells.Item(5, 21).Value = 1.51  # U5: 1.01 -> 1.51
$ws.Cells.Item(5, 22).Value = 1.07  # V5: 1.01 -> 1.07
$ws.Cells.Item(5, 23).Value = 2.68  # W5: 2.62 -> 2.68
$ws.Cells.Item(6, 14).Value = 4.1  # N6: 4 -> 4.1
$ws.Cells.Item(6, 15).Value = 1.28  # O6: 1.27 -> 1.28
$ws.Cells.Item(6, 16).Value = 2.08  # P6: 2.04 -> 2.08
$ws.Cells.Item(6, 18).Value = 1.43  # R6: 1.41 -> 1.43
$ws.Cells.Item(6, 20).Value = 1.67  # T6: 1.65 -> 1.67
$ws.Cells.Item(6, 21).Value = 2.24  # U6: 2.26 -> 2.24
$ws.Cells.Item(6, 23).Value = 1.7  # W6: 1.71 -> 1.7
$ws.Cells.Item(6, 25).Value = 14.5  # Y6: 15 -> 14.5
$ws.Cells.Item(6, 34).Value = 16.5  # AH6: 17 -> 16.5
$ws.Cells.Item(6, 37).Value = 24  # AK6: 25 -> 24
$ws.Cells.Item(6, 39).Value = 80  # AM6: 85 -> 80
$ws.Cells.Item(6, 40).Value = 17.5  # AN6: 18 -> 17.5
$ws.Cells.Item(7, 8).Value = 13  # H7: 12.5 -> 13
$ws.Cells.Item(7, 9).Value = 15  # I7: 16 -> 15
$ws.Cells.Item(7, 10).Value = 6.4  # J7: 6.2 -> 6.4
$ws.Cells.Item(7, 11).Value = 6.8  # K7: 7 -> 6.8
$ws.Cells.Item(7, 16).Value = 2.76  # P7: 2.78 -> 2.76
$ws.Cells.Item(7, 17).Value = 1.56  # Q7: 1.54 -> 1.56
$ws.Cells.Item(7, 18).Value = 1.69  # R7: 1.71 -> 1.69
$ws.Cells.Item(7, 19).Value = 2.3  # S7: 2.28 -> 2.3
$ws.Cells.Item(7, 20).Value = 2.2  # T7: 2.18 -> 2.2
$ws.Cells.Item(7, 21).Value = 1.75  # U7: 1.74 -> 1.75
$ws.Cells.Item(7, 25).Value = 50  # Y7: 55 -> 50
$ws.Cells.Item(7, 31).Value = 280  # AE7: 290 -> 280
$ws.Cells.Item(7, 32).Value = 8.4  # AF7: 8.6 -> 8.4
$ws.Cells.Item(7, 35).Value = 210  # AI7: 200 -> 210
$ws.Cells.Item(7, 36).Value = 9.4  # AJ7: 9.6 -> 9.4
$ws.Cells.Item(7, 39).Value = 210  # AM7: 200 -> 210
$ws.Cells.Item(7, 40).Value = 3.9  # AN7: 3.85 -> 3.9
$ws.Cells.Item(7, 41).Value = 390  # AO7: 1000 -> 390
$ws.Cells.Item(8, 14).Value = 2.92  # N8: 2.88 -> 2.92
$ws.Cells.Item(8, 16).Value = 1.79  # P8: 1.75 -> 1.79
$ws.Cells.Item(8, 18).Value = 1.28  # R8: 1.27 -> 1.28
$ws.Cells.Item(8, 24).Value = 13  # X8: 14 -> 13
$ws.Cells.Item(8, 27).Value = 220  # AA8: 240 -> 220
$ws.Cells.Item(8, 30).Value = 26  # AD8: 27 -> 26
$ws.Cells.Item(8, 31).Value = 120  # AE8: 130 -> 120
$ws.Cells.Item(8, 33).Value = 10.5  # AG8: 11.5 -> 10.5
$ws.Cells.Item(8, 37).Value = 1000  # AK8: 20 -> 1000
$ws.Cells.Item(8, 38).Value = 1000  # AL8: 50 -> 1000
$ws.Cells.Item(8, 40).Value = 12.5  # AN8: 1000 -> 12.5
$ws.Cells.Item(9, 14).Value = 3.7  # N9: 3.65 -> 3.7
$ws.Cells.Item(9, 18).Value = 1.36  # R9: 1.35 -> 1.36
$ws.Cells.Item(9, 20).Value = 1.82  # T9: 1.81 -> 1.82
$ws.Cells.Item(9, 21).Value = 2.06  # U9: 2.04 -> 2.06
$ws.Cells.Item(9, 27).Value = 100  # AA9: 110 -> 100
$ws.Cells.Item(9, 34).Value = 1000  # AH9: 20 -> 1000
$ws.Cells.Item(9, 39).Value = 110  # AM9: 120 -> 110
$ws.Cells.Item(10, 6).Value = 1.3  # F10: 1.32 -> 1.3
$ws.Cells.Item(10, 9).Value = 14  # I10: 14.5 -> 14
$ws.Cells.Item(10, 10).Value = 6  # J10: 5.9 -> 6
$ws.Cells.Item(10, 11).Value = 6.6  # K10: 6.8 -> 6.6
$ws.Cells.Item(10, 12).Value = 1.3  # L10: 1.01 -> 1.3
$ws.Cells.Item(10, 15).Value = 1.22  # O10: 1.21 -> 1.22
$ws.Cells.Item(10, 16).Value = 2.34  # P10: 2.38 -> 2.34
$ws.Cells.Item(10, 17).Value = 1.66  # Q10: 1.64 -> 1.66
$ws.Cells.Item(10, 19).Value = 2.62  # S10: 2.6 -> 2.62
$ws.Cells.Item(10, 20).Value = 2.2  # T10: 2.18 -> 2.2
$ws.Cells.Item(10, 28).Value = 9.2  # AB10: 9.6 -> 9.199999999999999
$ws.Cells.Item(10, 40).Value = 5.4  # AN10: 5.5 -> 5.4
$ws.Cells.Item(11, 6).Value = 1.89  # F11: 1.83 -> 1.89
$ws.Cells.Item(11, 9).Value = 4.6  # I11: 4.7 -> 4.6
$ws.Cells.Item(11, 11).Value = 4.2  # K11: 4.4 -> 4.2
$ws.Cells.Item(11, 12).Value = 1.5  # L11: 1.01 -> 1.5
$ws.Cells.Item(11, 14).Value = 3.05  # N11: 3 -> 3.05
$ws.Cells.Item(11, 19).Value = 4.2  # S11: 4.1 -> 4.2
$ws.Cells.Item(11, 24).Value = 12  # X11: 12.5 -> 12
$ws.Cells.Item(11, 25).Value = 1000  # Y11: 14.5 -> 1000
$ws.Cells.Item(11, 26).Value = 32  # Z11: 34 -> 32
$ws.Cells.Item(11, 28).Value = 1000  # AB11: 8.800000000000001 -> 1000
$ws.Cells.Item(11, 29).Value = 8.8  # AC11: 9.199999999999999 -> 8.800000000000001
$ws.Cells.Item(11, 30).Value = 19.5  # AD11: 21 -> 19.5
$ws.Cells.Item(12, 15).Value = 1.38  # O12: 1.37 -> 1.38
$ws.Cells.Item(12, 17).Value = 2  # Q12: 1.98 -> 2
$ws.Cells.Item(12, 19).Value = 3.9  # S12: 3.65 -> 3.9
$ws.Cells.Item(12, 20).Value = 1.84  # T12: 1.69 -> 1.84
$ws.Cells.Item(13, 8).Value = 2.68  # H13: 2.72 -> 2.68
$ws.Cells.Item(13, 15).Value = 1.48  # O13: 1.47 -> 1.48
$ws.Cells.Item(13, 17).Value = 1.48  # Q13: 1.47 -> 1.48
$ws.Cells.Item(13, 19).Value = 2.44  # S13: 3.9 -> 2.44
$ws.Cells.Item(13, 25).Value = 10.5  # Y13: 10 -> 10.5
$ws.Cells.Item(13, 26).Value = 18  # Z13: 17 -> 18
$ws.Cells.Item(14, 9).Value = 17.5  # I14: 18 -> 17.5
$ws.Cells.Item(14, 11).Value = 7.8  # K14: 8.199999999999999 -> 7.8
$ws.Cells.Item(14, 12).Value = 1.27  # L14: 1.01 -> 1.27
$ws.Cells.Item(14, 14).Value = 6.2  # N14: 5.9 -> 6.2
$ws.Cells.Item(14, 16).Value = 2.7  # P14: 2.62 -> 2.7
$ws.Cells.Item(14, 17).Value = 1.52  # Q14: 1.51 -> 1.52
$ws.Cells.Item(14, 19).Value = 2.32  # S14: 2.3 -> 2.32
$ws.Cells.Item(14, 21).Value = 1.72  # U14: 1.7 -> 1.72
$ws.Cells.Item(14, 24).Value = 1000  # X14: 30 -> 1000
$ws.Cells.Item(14, 27).Value = 970  # AA14: 990 -> 970
$ws.Cells.Item(14, 28).Value = 14.5  # AB14: 1000 -> 14.5
$ws.Cells.Item(14, 30).Value = 980  # AD14: 1000 -> 980
$ws.Cells.Item(14, 31).Value = 300  # AE14: 310 -> 300
$ws.Cells.Item(14, 35).Value = 220  # AI14: 1000 -> 220
$ws.Cells.Item(14, 36).Value = 9.2  # AJ14: 9.4 -> 9.199999999999999
$ws.Cells.Item(14, 37).Value = 13.5  # AK14: 17.5 -> 13.5
$ws.Cells.Item(14, 39).Value = 220  # AM14: 240 -> 220
$ws.Cells.Item(14, 40).Value = 4.2  # AN14: 4.1 -> 4.2
$ws.Cells.Item(14, 41).Value = 360  # AO14: 410 -> 360
$ws.Cells.Item(15, 6).Value = 1.83  # F15: 1.8 -> 1.83
$ws.Cells.Item(15, 7).Value = 1.95  # G15: 1.96 -> 1.95
$ws.Cells.Item(15, 9).Value = 5  # I15: 5.3 -> 5
$ws.Cells.Item(15, 10).Value = 3.65  # J15: 3.6 -> 3.65
$ws.Cells.Item(15, 11).Value = 4.1  # K15: 4.2 -> 4.1
$ws.Cells.Item(15, 14).Value = 3.25  # N15: 3.15 -> 3.25
$ws.Cells.Item(15, 15).Value = 1.36  # O15: 1.29 -> 1.36
$ws.Cells.Item(15, 16).Value = 1.77  # P15: 1.76 -> 1.77
$ws.Cells.Item(15, 17).Value = 1.95  # Q15: 1.9 -> 1.95
$ws.Cells.Item(15, 19).Value = 3.4  # S15: 3.35 -> 3.4
$ws.Cells.Item(15, 20).Value = 1.93  # T15: 1.92 -> 1.93
$ws.Cells.Item(15, 21).Value = 1.89  # U15: 1.8 -> 1.89
$ws.Cells.Item(15, 22).Value = 1.25  # V15: 1.23 -> 1.25
$ws.Cells.Item(15, 28).Value = 9.8  # AB15: 10 -> 9.800000000000001
$ws.Cells.Item(16, 10).Value = 3.55  # J16: 3.6 -> 3.55
$ws.Cells.Item(16, 12).Value = 1.43  # L16: 1.01 -> 1.43
$ws.Cells.Item(16, 14).Value = 3.75  # N16: 3.65 -> 3.75
$ws.Cells.Item(16, 16).Value = 1.92  # P16: 1.9 -> 1.92
$ws.Cells.Item(16, 21).Value = 2.14  # U16: 2.1 -> 2.14
$ws.Cells.Item(16, 25).Value = 14  # Y16: 14.5 -> 14
$ws.Cells.Item(16, 34).Value = 18.5  # AH16: 18 -> 18.5
$ws.Cells.Item(16, 38).Value = 40  # AL16: 38 -> 40
$ws.Cells.Item(16, 39).Value = 100  # AM16: 110 -> 100
$ws.Cells.Item(16, 41).Value = 46  # AO16: 42 -> 46
$ws.Cells.Item(17, 7).Value = 1.5  # G17: 1.52 -> 1.5
$ws.Cells.Item(17, 8).Value = 8.2  # H17: 8 -> 8.199999999999999
$ws.Cells.Item(17, 10).Value = 4.5  # J17: 4.4 -> 4.5
$ws.Cells.Item(17, 12).Value = 1.41  # L17: 1.01 -> 1.41
$ws.Cells.Item(17, 14).Value = 3.95  # N17: 3.9 -> 3.95
$ws.Cells.Item(17, 16).Value = 2.02  # P17: 1.98 -> 2.02
$ws.Cells.Item(17, 23).Value = 2.96  # W17: 2.92 -> 2.96
$ws.Cells.Item(17, 27).Value = 380  # AA17: 400 -> 380
$ws.Cells.Item(17, 31).Value = 180  # AE17: 190 -> 180
$ws.Cells.Item(17, 34).Value = 29  # AH17: 32 -> 29
$ws.Cells.Item(17, 37).Value = 19.5  # AK17: 20 -> 19.5
$ws.Cells.Item(17, 39).Value = 210  # AM17: 220 -> 210
$ws.Cells.Item(17, 40).Value = 7.8  # AN17: 8 -> 7.8
$ws.Cells.Item(18, 10).Value = 8.2  # J18: 8 -> 8.199999999999999
$ws.Cells.Item(18, 17).Value = 1.52  # Q18: 1.5 -> 1.52
$ws.Cells.Item(18, 19).Value = 2.24  # S18: 2.22 -> 2.24
$ws.Cells.Item(18, 20).Value = 2.42  # T18: 2.4 -> 2.42
$ws.Cells.Item(18, 26).Value = 300  # Z18: 310 -> 300
$ws.Cells.Item(19, 10).Value = 3.7  # J19: 3.65 -> 3.7
$ws.Cells.Item(19, 12).Value = 1.37  # L19: 1.38 -> 1.37
$ws.Cells.Item(19, 24).Value = 16  # X19: 16.5 -> 16
$ws.Cells.Item(19, 28).Value = 15.5  # AB19: 16 -> 15.5
$ws.Cells.Item(19, 39).Value = 80  # AM19: 85 -> 80
$ws.Cells.Item(19, 41).Value = 15.5  # AO19: 16 -> 15.5
$ws.Cells.Item(20, 14).Value = 7  # N20: 6.8 -> 7
$ws.Cells.Item(20, 16).Value = 3.05  # P20: 2.94 -> 3.05
$ws.Cells.Item(20, 18).Value = 1.82  # R20: 1.79 -> 1.82
$ws.Cells.Item(20, 20).Value = 1.85  # T20: 1.84 -> 1.85
$ws.Cells.Item(20, 25).Value = 60  # Y20: 1000 -> 60
$ws.Cells.Item(20, 28).Value = 14.5  # AB20: 15.5 -> 14.5
$ws.Cells.Item(20, 29).Value = 15.5  # AC20: 16 -> 15.5
$ws.Cells.Item(20, 31).Value = 190  # AE20: 200 -> 190
$ws.Cells.Item(20, 34).Value = 27  # AH20: 28 -> 27
$ws.Cells.Item(20, 35).Value = 120  # AI20: 130 -> 120
$ws.Cells.Item(20, 36).Value = 13.5  # AJ20: 11.5 -> 13.5
$ws.Cells.Item(20, 37).Value = 13  # AK20: 14 -> 13
$ws.Cells.Item(20, 39).Value = 120  # AM20: 140 -> 120
$ws.Cells.Item(20, 40).Value = 3.75  # AN20: 3.85 -> 3.75
$ws.Cells.Item(21, 7).Value = 1.63  # G21: 1.64 -> 1.63
$ws.Cells.Item(21, 9).Value = 6.4  # I21: 6.6 -> 6.4
$ws.Cells.Item(21, 10).Value = 4.5  # J21: 4.4 -> 4.5
$ws.Cells.Item(21, 16).Value = 2.44  # P21: 2.4 -> 2.44
$ws.Cells.Item(21, 17).Value = 1.61  # Q21: 1.59 -> 1.61
$ws.Cells.Item(21, 19).Value = 2.5  # S21: 2.46 -> 2.5
$ws.Cells.Item(21, 20).Value = 1.73  # T21: 1.7 -> 1.73
$ws.Cells.Item(21, 21).Value = 2.2  # U21: 2.18 -> 2.2
$ws.Cells.Item(21, 23).Value = 2.58  # W21: 2.56 -> 2.58
$ws.Cells.Item(21, 31).Value = 90  # AE21: 75 -> 90
$ws.Cells.Item(21, 35).Value = 65  # AI21: 1000 -> 65
$ws.Cells.Item(21, 40).Value = 6.6  # AN21: 6.8 -> 6.6
$ws.Cells.Item(22, 6).Value = 1.24  # F22: 1.23 -> 1.24
$ws.Cells.Item(22, 7).Value = 1.25  # G22: 1.26 -> 1.25
$ws.Cells.Item(22, 10).Value = 7.6  # J22: 7.4 -> 7.6
$ws.Cells.Item(22, 11).Value = 8  # K22: 8.4 -> 8
$ws.Cells.Item(22, 14).Value = 8.6  # N22: 7.8 -> 8.6
$ws.Cells.Item(22, 15).Value = 1.11  # O22: 1.12 -> 1.11
$ws.Cells.Item(22, 16).Value = 3.35  # P22: 3.3 -> 3.35
$ws.Cells.Item(22, 18).Value = 1.94  # R22: 1.97 -> 1.94
$ws.Cells.Item(22, 19).Value = 1.92  # S22: 1.89 -> 1.92
$ws.Cells.Item(22, 20).Value = 1.89  # T22: 1.88 -> 1.89
$ws.Cells.Item(22, 28).Value = 16  # AB22: 18.5 -> 16
$ws.Cells.Item(22, 32).Value = 11  # AF22: 10.5 -> 11
$ws.Cells.Item(22, 39).Value = 140  # AM22: 150 -> 140
$ws.Cells.Item(22, 40).Value = 3.6  # AN22: 3.55 -> 3.6
$ws.Cells.Item(23, 6).Value = 1.65  # F23: 1.64 -> 1.65
$ws.Cells.Item(23, 9).Value = 5.8  # I23: 6 -> 5.8
$ws.Cells.Item(23, 10).Value = 4.1  # J23: 4.2 -> 4.1
$ws.Cells.Item(23, 15).Value = 1.24  # O23: 1.23 -> 1.24
$ws.Cells.Item(23, 16).Value = 2.22  # P23: 2.2 -> 2.22
$ws.Cells.Item(23, 19).Value = 2.72  # S23: 2.66 -> 2.72
$ws.Cells.Item(23, 20).Value = 1.74  # T23: 1.72 -> 1.74
$ws.Cells.Item(23, 21).Value = 2.16  # U23: 2.12 -> 2.16
$ws.Cells.Item(23, 23).Value = 2.32  # W23: 2.34 -> 2.32
$ws.Cells.Item(23, 29).Value = 10  # AC23: 1000 -> 10
$ws.Cells.Item(23, 34).Value = 22  # AH23: 1000 -> 22
